$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $style = $cell.Style
    if ($value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $cell.Value = "'" + $value
    } else {
        $cell.Value = $value
    }
    $cell.Style = $style
}

# --- Price/Volume updates for unchanged coin rows ---
Set-TextValue $ws.Range("D2") "28.277.02"
Set-TextValue $ws.Range("E2") "  -1.45%  "
Set-TextValue $ws.Range("D3") "1.824.33"
Set-TextValue $ws.Range("E3") "  +0.35%  "
Set-TextValue $ws.Range("D4") "0.9952"
Set-TextValue $ws.Range("E4") "  -0.56%  "
Set-TextValue $ws.Range("D5") "325.93"
Set-TextValue $ws.Range("E5") "  -1.10%  "
Set-TextValue $ws.Range("D6") "0.9905"
Set-TextValue $ws.Range("E6") "  -0.77%  "
Set-TextValue $ws.Range("D7") "0.4436"
Set-TextValue $ws.Range("E7") "  -0.17%  "
Set-TextValue $ws.Range("D8") "0.3777"
Set-TextValue $ws.Range("E8") "  -1.25%  "
Set-TextValue $ws.Range("D9") "45.29"
Set-TextValue $ws.Range("E9") "  +0.80%  "
Set-TextValue $ws.Range("D10") "0.07749"
Set-TextValue $ws.Range("E10") "  +0.81%  "
Set-TextValue $ws.Range("E11") "  -1.75%  "
Set-TextValue $ws.Range("D12") "22.14"
Set-TextValue $ws.Range("E12") "  -3.73%  "
Set-TextValue $ws.Range("D13") "0.9886"
Set-TextValue $ws.Range("E13") "  -1.08%  "
Set-TextValue $ws.Range("D14") "6.283"
Set-TextValue $ws.Range("E14") "  -1.55%  "
Set-TextValue $ws.Range("D15") "7.515"
Set-TextValue $ws.Range("E15") "  -1.67%  "
Set-TextValue $ws.Range("D16") "1.813.04"
Set-TextValue $ws.Range("E16") "  +0.10%  "
Set-TextValue $ws.Range("D17") "92.26"
Set-TextValue $ws.Range("E17") "  +12.84%  "
Set-TextValue $ws.Range("E18") "  -1.64%  "
Set-TextValue $ws.Range("D19") "0.06360"
Set-TextValue $ws.Range("E19") "  -6.14%  "
Set-TextValue $ws.Range("D20") "0.9928"
Set-TextValue $ws.Range("E20") "  -0.60%  "
Set-TextValue $ws.Range("D21") "17.55"
Set-TextValue $ws.Range("E21") "  -2.12%  "
Set-TextValue $ws.Range("D22") "6.346"
Set-TextValue $ws.Range("E22") "  -0.66%  "
Set-TextValue $ws.Range("D23") "0.5331"
Set-TextValue $ws.Range("E23") "  -2.11%  "
Set-TextValue $ws.Range("D24") "28.336.32"
Set-TextValue $ws.Range("E24") "  -1.23%  "
Set-TextValue $ws.Range("D25") "11.67"
Set-TextValue $ws.Range("E25") "  -2.59%  "
Set-TextValue $ws.Range("D26") "2.118"
Set-TextValue $ws.Range("E26") "  -12.46%  "
Set-TextValue $ws.Range("D27") "20.91"
Set-TextValue $ws.Range("E27") "  +0.05%  "
Set-TextValue $ws.Range("D28") "153.40"
Set-TextValue $ws.Range("E28") "  -0.19%  "
Set-TextValue $ws.Range("D29") "2.373"
Set-TextValue $ws.Range("E29") "  -1.15%  "
Set-TextValue $ws.Range("D30") "2.015.66"
Set-TextValue $ws.Range("E30") "  -0.06%  "
Set-TextValue $ws.Range("D31") "129.73"
Set-TextValue $ws.Range("E31") "  -2.99%  "
Set-TextValue $ws.Range("D32") "1.211"
Set-TextValue $ws.Range("E32") "  -6.93%  "
Set-TextValue $ws.Range("D33") "5.858"
Set-TextValue $ws.Range("E33") "  -1.15%  "
Set-TextValue $ws.Range("D34") "0.09236"
Set-TextValue $ws.Range("E34") "  -1.30%  "
Set-TextValue $ws.Range("D35") "3.644"
Set-TextValue $ws.Range("E35") "  -8.14%  "
Set-TextValue $ws.Range("D36") "12.81"
Set-TextValue $ws.Range("E36") "  +3.33%  "
Set-TextValue $ws.Range("D37") "0.02348"
Set-TextValue $ws.Range("E37") "  -0.71%  "
Set-TextValue $ws.Range("E38") "  -4.66%  "
Set-TextValue $ws.Range("D41") "0.06208"
Set-TextValue $ws.Range("E41") "  -3.52%  "
Set-TextValue $ws.Range("D42") "1.186"
Set-TextValue $ws.Range("E42") "  -2.47%  "
Set-TextValue $ws.Range("D43") "8.053"
Set-TextValue $ws.Range("E43") "  -2.35%  "
Set-TextValue $ws.Range("D46") "13.94"
Set-TextValue $ws.Range("E46") "  -2.02%  "
Set-TextValue $ws.Range("D47") "0.6100"
Set-TextValue $ws.Range("E47") "  -1.50%  "
Set-TextValue $ws.Range("D48") "3.735"
Set-TextValue $ws.Range("D49") "126.90"
Set-TextValue $ws.Range("E49") "  -2.31%  "
Set-TextValue $ws.Range("D50") "2.029"

# --- Rows with coin swaps / replacements (B, C, D, E) ---
Set-TextValue $ws.Range("B39") "TheSandbox"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D39") "0.6591"
Set-TextValue $ws.Range("E39") "  -2.31%  "
Set-TextValue $ws.Range("B40") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D40") "5.168"
Set-TextValue $ws.Range("E40") "  -2.10%  "
Set-TextValue $ws.Range("B44") "WEMIXTOKEN"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "1.405"
Set-TextValue $ws.Range("E44") "  -3.36%  "
Set-TextValue $ws.Range("B45") "Frax"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D45") "0.9891"
Set-TextValue $ws.Range("E45") "  -0.93%  "
Set-TextValue $ws.Range("B51") "Cronos"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.07024"
Set-TextValue $ws.Range("E51") "  -1.57%  "
